# The commit adds a "Date and Time" row at the very top and a
# "Cycle Count of battery" row further down, which renumbers every row
# below each insertion point; the underlying metric values were also
# recalculated upstream (not just shifted). Rather than inserting rows and
# hoping formats/values line up, write the full resulting A1:B45 block
# directly, row by row, from the known target state.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Date and Time"
$ws.Cells.Item(1, 2).ClearFormats()
$ws.Cells.Item(1, 2).Value = "2024-03-11 16:05:22.553000 to 2024-03-11 16:57:30.743000"

$ws.Cells.Item(2, 1).Value = "Total time taken for the ride"
$ws.Cells.Item(2, 2).Value = 0.03623425925925926
$ws.Cells.Item(2, 2).NumberFormat = "[hh]:mm:ss"

$ws.Cells.Item(3, 1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3, 2).Value = 32.97601722222222

$ws.Cells.Item(4, 1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4, 2).Value = 1673.844472675

$ws.Cells.Item(5, 1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5, 2).Value = 39.268

$ws.Cells.Item(6, 1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6, 2).Value = 7.33

$ws.Cells.Item(7, 1).Value = "Starting SoC (%)"
$ws.Cells.Item(7, 2).Value = 18

$ws.Cells.Item(8, 1).Value = "Ending SoC (%)"
$ws.Cells.Item(8, 2).Value = 99

$ws.Cells.Item(9, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(9, 2).Value = 30.98661569075563

$ws.Cells.Item(10, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 2).Value = 54.01830549614895

$ws.Cells.Item(11, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11, 2).Value = 81

$ws.Cells.Item(12, 1).Value = "Mode"
$ws.Cells.Item(12, 2).Value = "Custom mode`n69.63%`nEco mode`n18.00%`nSports mode`n0.06%"

$ws.Cells.Item(13, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 2).Value = 5421.410500000001

$ws.Cells.Item(14, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 2).Value = -1935.080315231214

$ws.Cells.Item(15, 1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15, 2).Value = 1.447407298055555

$ws.Cells.Item(16, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16, 2).Value = 0.08639732069129558

$ws.Cells.Item(17, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.332

$ws.Cells.Item(18, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18, 2).Value = 3.071

$ws.Cells.Item(19, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 2).Value = 0.2609999999999997

$ws.Cells.Item(20, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 2).Value = 38

$ws.Cells.Item(21, 1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21, 2).Value = 48

$ws.Cells.Item(22, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22, 2).Value = 10

$ws.Cells.Item(23, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 2).Value = 70

$ws.Cells.Item(24, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 2).Value = 66

$ws.Cells.Item(25, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 2).Value = 65

$ws.Cells.Item(26, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 2).Value = 68

$ws.Cells.Item(27, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 2).Value = 100

$ws.Cells.Item(28, 1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28, 2).Value = 0

$ws.Cells.Item(29, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 2).Value = 48

$ws.Cells.Item(30, 1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30, 2).Value = 38

$ws.Cells.Item(31, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31, 2).Value = 10

$ws.Cells.Item(32, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32, 2).Value = 53

$ws.Cells.Item(33, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33, 2).Value = 1.747728912777778

$ws.Cells.Item(34, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34, 2).Value = [double]"1.551055123160967e-07"

$ws.Cells.Item(35, 1).Value = "Cycle Count of battery"
$ws.Cells.Item(35, 2).Value = 27

$ws.Cells.Item(36, 1).Value = "Idling time percentage"
$ws.Cells.Item(36, 2).Value = 15.93362445414847

$ws.Cells.Item(37, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37, 2).Value = 8.419213973799126

$ws.Cells.Item(38, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38, 2).Value = 3.814847161572053

$ws.Cells.Item(39, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39, 2).Value = 8.978165938864628

$ws.Cells.Item(40, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40, 2).Value = 17.75371179039301

$ws.Cells.Item(41, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41, 2).Value = 11.83231441048035

$ws.Cells.Item(42, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42, 2).Value = 11.7414847161572

$ws.Cells.Item(43, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43, 2).Value = 13.95633187772926

$ws.Cells.Item(44, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44, 2).Value = 7.479475982532752

$ws.Cells.Item(45, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45, 2).Value = 0
